# Add 3 new data sheets (Total_Population, Race_Ethnicity, Age) to the
# scag_tabbed workbook, and register them on the TOC sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# TOC sheet: append rows describing the three new sheets
# ---------------------------------------------------------------------------
$toc = $wb.Worksheets.Item("TOC")

$toc.Range("A13").Value = "Total_Population"
$toc.Range("B13").Value = "Total Population by County and SCAG Region"

$toc.Range("A14").Value = "Race_Ethnicity"
$toc.Range("B14").Value = "Race/Ethnicity Distribution (%) by County and SCAG Region"

$toc.Range("A15").Value = "Age"
$toc.Range("B15").Value = "Age Distribution (%) by County and SCAG Region"

# ---------------------------------------------------------------------------
# helper to bold + center a header row range
# ---------------------------------------------------------------------------
function Format-Header($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108  # xlCenter
}

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------------
# Sheet 13: Total_Population
# ---------------------------------------------------------------------------
$popSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$popSheet.Name = "Total_Population"

$popSheet.Range("A1").Value = "county"
$popSheet.Range("B1").Value = "total_pop"
Format-Header $popSheet.Range("A1:B1")

$popData = @(
    @("Imperial", 179943),
    @("Los Angeles", 10019738),
    @("Orange", 3182954),
    @("Riverside", 2409370),
    @("San Bernardino", 2170489),
    @("Ventura", 844838),
    @("SCAG", 18807332)
)

for ($i = 0; $i -lt $popData.Count; $i++) {
    $row = $i + 2
    $popSheet.Cells.Item($row, 1).Value = $popData[$i][0]
    $popSheet.Cells.Item($row, 2).Value = $popData[$i][1]
}

# ---------------------------------------------------------------------------
# Sheet 14: Race_Ethnicity
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$raceSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$raceSheet.Name = "Race_Ethnicity"

$raceHeaders = @("race", "Imperial", "Los Angeles", "Orange", "Riverside", "San Bernardino", "Ventura", "SCAG")
for ($c = 0; $c -lt $raceHeaders.Count; $c++) {
    $raceSheet.Cells.Item(1, $c + 1).Value = $raceHeaders[$c]
}
Format-Header $raceSheet.Range("A1:H1")

$raceData = @(
    @("Asian/Pacific Islander", 1.33, 14.76, 21.37, 6.77, 7.56, 7.32, 13.56),
    @("Black", 2.47, 7.69, 1.57, 6.18, 7.61, 1.7, 6.13),
    @("Hispanic/Latino", 85.04000000000001, 48.7, 34.02, 50.28, 54.59, 43.3, 47.2),
    @("Multiracial/Other", 0.66, 3.18, 3.86, 3.22, 3.31, 3.29, 3.29),
    @("Native American", 0.6899999999999999, 0.2, 0.15, 0.35, 0.31, 0.22, 0.23),
    @("White", 9.800000000000001, 25.47, 39.03, 33.21, 26.62, 44.17, 29.58)
)

for ($i = 0; $i -lt $raceData.Count; $i++) {
    $row = $i + 2
    for ($c = 0; $c -lt $raceData[$i].Count; $c++) {
        $raceSheet.Cells.Item($row, $c + 1).Value = $raceData[$i][$c]
    }
}

# ---------------------------------------------------------------------------
# Sheet 15: Age
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ageSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ageSheet.Name = "Age"

$ageSheet.Range("A1").Value = "age_categories"
$ageSheet.Range("B1").Value = "county"
$ageSheet.Range("C1").Value = "percentages"
Format-Header $ageSheet.Range("A1:C1")

$ageData = @(
    @("<18 years", "Imperial", 28.7),
    @("<18 years", "Los Angeles", 21.64),
    @("<18 years", "Orange", 21.97),
    @("<18 years", "Riverside", 25.16),
    @("<18 years", "San Bernardino", 26.43),
    @("<18 years", "Ventura", 22.87),
    @("<18 years", "SCAG", 22.82),
    @("18 - 64 years", "Imperial", 58.46),
    @("18 - 64 years", "Los Angeles", 64.66),
    @("18 - 64 years", "Orange", 63.18),
    @("18 - 64 years", "Riverside", 60.49),
    @("18 - 64 years", "San Bernardino", 61.98),
    @("18 - 64 years", "Ventura", 61.48),
    @("18 - 64 years", "SCAG", 63.36),
    @("65+ years", "Imperial", 12.84),
    @("65+ years", "Los Angeles", 13.71),
    @("65+ years", "Orange", 14.85),
    @("65+ years", "Riverside", 14.35),
    @("65+ years", "San Bernardino", 11.59),
    @("65+ years", "Ventura", 15.65),
    @("65+ years", "SCAG", 13.82)
)

for ($i = 0; $i -lt $ageData.Count; $i++) {
    $row = $i + 2
    $ageSheet.Cells.Item($row, 1).Value = $ageData[$i][0]
    $ageSheet.Cells.Item($row, 2).Value = $ageData[$i][1]
    $ageSheet.Cells.Item($row, 3).Value = $ageData[$i][2]
}

# Leave the original active sheet/selection as it was (TOC, A1 selected)
$toc.Activate()
$toc.Range("A1").Select()
